$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.148.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.56%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.430.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.83%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.77%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.76%  '

# Row 7
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.427.70'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.84%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.120'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.76%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.06'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.48%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.373'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.55%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.018.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000178'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -9.43%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.02%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.440.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.45%  '

# Row 17
$ws.Range("E17").Value = '  -2.06%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.187.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -12.57%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.91%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.17%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.01%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.550'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.91%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.52%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.573.35'
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = '  -7.93%  '

# Row 28
$ws.Range("E28").Value = '  -0.02%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.62%  '

# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.25%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.442.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.01%  '

# Row 33
$ws.Range("E33").Value = '  +0.03%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.143'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.94%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.29%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '169.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.44%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.74'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.55%  '

# Row 38
$ws.Range("E38").Value = '  -11.60%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.28%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.36%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0753'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.44%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.809'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.52%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.89%  '

# Row 44
$ws.Range("E44").Value = '  -0.22%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -14.74%  '

# Row 46
$ws.Range("E46").Value = '  -9.00%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.41%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.17%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.45'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.92%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -13.49%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.161.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.13%  '
